$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "What is something you would like to see to make a past experience
#    better than your first time playing?"
#    -> "... make an experience better ..."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("a past experience", $true, $false, $false, $false, `
    $false, $true, 1, $false, "an experience", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Move the hidden "_GoBack" bookmark from the paragraph that currently
#    holds it ("... quick tutorial ...") down to the paragraph that will
#    become the new final paragraph of the body ("... adult content ...").
#    This mirrors the net movement shown in the diff once the duplicated
#    Q/A block + "Core Feature" list (which originally followed it) are
#    removed below.
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "*sectored for adults only*") {
        $targetPara = $d.Paragraphs($i)
        break
    }
}
# A zero-length Range confuses this host's Bookmarks.Add (it silently
# anchors at document position 0), so a collapsed insertion point can't be
# bookmarked directly. Work around it by inserting a 1-character sentinel
# at the true end of the paragraph's text, bookmarking that single
# character, then deleting the sentinel again -- the bookmark collapses
# back down to an empty range but stays anchored in the right spot, same
# as the bookmarkStart/bookmarkEnd pair sitting right after the last run.
$endPos = $targetPara.Range.End - 1
$sentinel = $d.Range($endPos, $endPos)
$sentinel.InsertBefore("\u0001")
$sentinelRange = $d.Range($endPos, $endPos + 1)
$d.Bookmarks.Add("_GoBack", $sentinelRange) | Out-Null
$d.Range($endPos, $endPos + 1).Delete()

# ---------------------------------------------------------------------------
# 3) Delete the trailing "Core Feature" block (header, blank lines and the
#    four bulleted list items) that followed the duplicated Q/A content.
#    That whole stretch, starting with the blank paragraph right after the
#    "... adult content ..." answer and running to the end of the body, is
#    removed.
# ---------------------------------------------------------------------------
$startPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "*Core Feature*") {
        $startPara = $d.Paragraphs($i - 1)
        break
    }
}
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$killRange = $d.Range($startPara.Range.Start, $lastPara.Range.End)
$killRange.Delete()
